$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Palabras")

$ws.Range("D2").Value = 9
$ws.Range("D3").Value = 22
$ws.Range("D6").Value = 3
$ws.Range("D7").Value = 3
$ws.Range("D8").Value = 22
$ws.Range("D9").Value = 1
$ws.Range("D11").Value = 11
$ws.Range("D12").Value = 4
$ws.Range("D13").Value = 3
$ws.Range("D14").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("D19").Value = 11
$ws.Range("D20").Value = 15
$ws.Range("D21").Value = 6
$ws.Range("D22").Value = 8
$ws.Range("D23").Value = 10
$ws.Range("D24").Value = 9
$ws.Range("D25").Value = 5
$ws.Range("D26").Value = 6
